$d = $word.ActiveDocument

# Find the paragraph that contains "Answer Scheme" (bold heading that starts
# the answer-scheme block we want to remove), and the last numbered answer
# paragraph ("10.    A"), then delete the whole range spanning them while
# leaving the preceding page-break paragraph and the trailing sectPr intact.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $text = $p.Range.Text
  if ($text -match "Answer Scheme") {
    $startPara = $p
  }
  if ($text -match "^10\.\s*A") {
    $endPara = $p
  }
}

if ($startPara -ne $null -and $endPara -ne $null) {
  $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
  $range.Delete()
}
